$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the scrolled/selected view to match the new state
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("C46").Select()

# C6: 0.8 -> 1 (percentage formatted cell keeps its style, only the value changes)
$ws.Range("C6").Value = 1

# C43 used to hold the text "en proceso" (shared string); it now becomes a
# numeric 1 formatted the same way as the other percentage cells (style s="3"),
# so copy the number format from a sibling cell that already has it.
$ws.Range("C44").Copy() | Out-Null
$ws.Range("C43").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C43").Value = 1
$excel.CutCopyMode = 0

# Row 45 gains a responsable ("Agustina") and a status ("en proceso"),
# matching the pattern used by the other rows in this table.
$ws.Range("B45").Value = "Agustina"
$ws.Range("C45").Value = "en proceso"
